$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4703.846
$ws.Range("J38").Value = 6088
$ws.Range("L38").Value = 18264
$ws.Range("N38").Value = -19008
$ws.Range("H39").Value = 299.75
$ws.Range("I39").Value = 98.28570999999999
$ws.Range("J39").Value = 581.8
$ws.Range("K39").Value = 294.85713
$ws.Range("L39").Value = 1745.4
$ws.Range("M39").Value = 1.142870000000016
$ws.Range("N39").Value = -2337.4
$ws.Range("H40").Value = 1800
$ws.Range("J40").Value = 1800
$ws.Range("L40").Value = 1800
$ws.Range("N40").Value = -2150
$ws.Range("H137").Value = 1906963.5
$ws.Range("I137").Value = 2268718.5
$ws.Range("K137").Value = 6806155.5
$ws.Range("M137").Value = -6803605.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 7000
$ws.Range("J56").Value = 8950
$ws.Range("L56").Value = 8950
$ws.Range("N56").Value = -10434
$ws.Range("H110").Value = 1488.2858
$ws.Range("I110").Value = 1488.2858
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1488.2858
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 556.7141999999999
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2517530.2
$ws.Range("I7").Value = 9645
$ws.Range("J7").Value = 2875799.5
$ws.Range("K7").Value = 9645
$ws.Range("L7").Value = 2875799.5
$ws.Range("M7").Value = -9532
$ws.Range("N7").Value = -2876025.5
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20454
$ws.Range("H75").Value = 6500
$ws.Range("I75").Value = 4750
$ws.Range("J75").Value = 10000
$ws.Range("K75").Value = 4750
$ws.Range("L75").Value = 10000
$ws.Range("M75").Value = -3814
$ws.Range("N75").Value = -11872
$ws.Range("H78").Value = 6500
$ws.Range("I78").Value = 4750
$ws.Range("J78").Value = 10000
$ws.Range("K78").Value = 14250
$ws.Range("L78").Value = 30000
$ws.Range("M78").Value = -9570
$ws.Range("N78").Value = -39360
$ws.Range("H82").Value = 30144.5
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 32308
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 32308
$ws.Range("M82").Value = -14617
$ws.Range("N82").Value = -33074
$ws.Range("H85").Value = 30144.5
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 32308
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 32308
$ws.Range("M85").Value = -13674
$ws.Range("N85").Value = -34960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 50000
$ws.Range("J48").Value = 50000
$ws.Range("L48").Value = 50000
$ws.Range("N48").Value = -50952
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H95").Value = 29562
$ws.Range("J95").Value = 29562
$ws.Range("L95").Value = 29562
$ws.Range("N95").Value = -35054
$ws.Range("H105").Value = 1620.0476
$ws.Range("I105").Value = 1278.9445
$ws.Range("J105").Value = 3666.6667
$ws.Range("K105").Value = 1278.9445
$ws.Range("L105").Value = 3666.6667
$ws.Range("M105").Value = 468.0554999999999
$ws.Range("N105").Value = -7160.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H101").Value = 2000
$ws.Range("J101").Value = 2000
$ws.Range("L101").Value = 6000
$ws.Range("N101").Value = -10868
$ws.Range("H121").Value = 1786.6232
$ws.Range("J121").Value = 1867.9539
$ws.Range("L121").Value = 5603.861699999999
$ws.Range("N121").Value = -8223.861699999999
$ws.Range("H125").Value = 4560
$ws.Range("I125").Value = 1620
$ws.Range("J125").Value = 7500
$ws.Range("K125").Value = 4860
$ws.Range("L125").Value = 22500
$ws.Range("M125").Value = 60
$ws.Range("N125").Value = -32340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20560
$ws.Range("H49").Value = 9979.429
$ws.Range("J49").Value = 12800.4
$ws.Range("L49").Value = 12800.4
$ws.Range("N49").Value = -13168.4
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -20996
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H120").Value = 38566.668
$ws.Range("J120").Value = 38566.668
$ws.Range("L120").Value = 38566.668
$ws.Range("N120").Value = -48242.668
$ws.Range("H126").Value = 3264.91
$ws.Range("I126").Value = 2905.0908
$ws.Range("J126").Value = 4469.522
$ws.Range("K126").Value = 8715.2724
$ws.Range("L126").Value = 13408.566
$ws.Range("M126").Value = -6245.2724
$ws.Range("N126").Value = -18348.566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5216.364
$ws.Range("I7").Value = 4295.8335
$ws.Range("J7").Value = 6321
$ws.Range("K7").Value = 4295.8335
$ws.Range("L7").Value = 6321
$ws.Range("M7").Value = -4183.8335
$ws.Range("N7").Value = -6545
$ws.Range("H126").Value = 5216.364
$ws.Range("I126").Value = 4295.8335
$ws.Range("J126").Value = 6321
$ws.Range("K126").Value = 12887.5005
$ws.Range("L126").Value = 18963
$ws.Range("M126").Value = -10417.5005
$ws.Range("N126").Value = -23903

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 32769
$ws.Range("J47").Value = 32769
$ws.Range("L47").Value = 32769
$ws.Range("N47").Value = -33913
$ws.Range("H122").Value = 6125.25
$ws.Range("I122").Value = 2750.5
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 8251.5
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -5801.5
$ws.Range("N122").Value = -33400
